$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    44 = -3.206799444181351
    45 = -3.209591218444118
    50 = -3.223550089757942
    51 = -3.226341864020707
    52 = -3.229133638283471
    53 = -3.217146336696169
    54 = -3.2017146768479
    81 = -3.346376550581282
    82 = -3.390412083964656
    83 = -3.400858494119914
    84 = -3.385058172577313
    85 = -3.358867812671553
    86 = -3.323150665129971
    87 = -3.254031234515312
    88 = -3.173450656128138
    89 = -3.079467799254104
    92 = -3.435121590036689
    93 = -3.476134744936997
    94 = -3.499937344937792
    95 = -3.507223093128693
    96 = -3.487569036546394
    97 = -3.450296227054411
    98 = -3.383770308076794
    99 = -3.292353145717887
    103 = -3.520009958178731
    104 = -3.555456197634888
    105 = -3.571213448440229
    106 = -3.563680012309287
    107 = -3.537671234341154
    108 = -3.492621859105395
    109 = -3.428483775593861
    110 = -3.350970152897967
    111 = -3.258498600038136
    114 = -3.561441234999713
    115 = -3.596457130478383
    116 = -3.61395424878778
    117 = -3.611581857620827
    118 = -3.588096510605501
    119 = -3.543053006215779
    120 = -3.479938474915214
    121 = -3.404747227695209
    122 = -3.309567466726047
    125 = -3.603633701674299
    126 = -3.633595315254041
    127 = -3.647786744280726
    128 = -3.639343491996382
    129 = -3.611648645276802
    130 = -3.566065272752695
    131 = -3.504169106570375
    132 = -3.43039418076205
    133 = -3.343200634906486
    136 = -3.628597995717426
    137 = -3.657244825031758
    138 = -3.671349806927445
    139 = -3.663960399829639
    140 = -3.636301172014854
    141 = -3.59129686709689
    142 = -3.530669665833261
    143 = -3.45862438636013
    144 = -3.372098767043833
    147 = -3.653522909296032
    148 = -3.678466205346917
    149 = -3.69003320976123
    150 = -3.679508644539837
    151 = -3.651026265702303
    152 = -3.605796704641904
    153 = -3.545831926400007
    154 = -3.474688083899361
    155 = -3.392899699789125
    158 = -3.669508263422132
    159 = -3.69360288522447
    160 = -3.704345297926025
    161 = -3.693724885187034
    162 = -3.666027597235246
    163 = -3.621863113637593
    164 = -3.561668729517023
    165 = -3.492373876331519
    166 = -3.411517128213375
    169 = -3.685335275268046
    170 = -3.706557895154997
    171 = -3.715521867279732
    172 = -3.703998200123845
    173 = -3.675450307186651
    174 = -3.630985534087508
    175 = -3.571991554490325
    176 = -3.503005852977148
    177 = -3.424859438971547
    180 = -3.69613571111312
    181 = -3.716567887516254
    182 = -3.724752463089372
    183 = -3.71294549342367
    184 = -3.68508408307939
    185 = -3.641237101576034
    186 = -3.582784309966356
    187 = -3.51497134444751
    188 = -3.437390249708752
    191 = -3.706975436092137
    192 = -3.725534547929304
    193 = -3.731886820719709
    194 = -3.719606922546133
    195 = -3.690988038200021
    196 = -3.648022230543787
    197 = -3.590175377571621
    198 = -3.522840734462941
    199 = -3.446176991636912
    202 = -3.714818669084385
    203 = -3.732692674570675
    204 = -3.738022796925615
    205 = -3.725696204525273
    206 = -3.697813399117424
    207 = -3.654923866186364
    208 = -3.598479676051884
    209 = -3.530987063775228
    210 = -3.455238735908447
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
